$d = $word.ActiveDocument

# The final paragraph of the document holds the "Break the Problem Apart" text
# for the "Predicting Fingers" problem. Its pPr carries a stray empty <w:rPr/>
# that the target revision drops, so we recreate the paragraph from scratch
# (delete + reinsert) rather than editing it in place, and then append the new
# "Identify Potential Solutions" paragraphs the same fresh way.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastText = $lastPara.Range.Text
$lastText = $lastText.TrimEnd([char]13)

$prevPara = $d.Paragraphs.Item($count - 1)
$prevEnd = $prevPara.Range.End

$fullRange = $d.Range($prevEnd, $lastPara.Range.End)
$fullRange.Delete()

$anchor = $d.Paragraphs.Item($d.Paragraphs.Count).Range

$newTexts = @(
    $lastText,
    "",
    "3) Identify Potential Solutions",
    "a) We could develop an equation that would tell us which figure any given number will fall on.",
    "b) We could manually count to each number using the same method as the girl.",
    "c) We could guess and hope for the best."
)

foreach ($t in $newTexts) {
    $anchor.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r = $newPara.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $t
    $anchor = $newPara.Range
}
